$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Cells in column D ("Price") store plain decimal-looking numbers
# (e.g. "322.87") as TEXT in the source data. Excel's default General-format
# input parsing would otherwise silently convert such text into a numeric
# value, so for every Price cell whose new value parses as a plain number we
# first force the cell to Text format ("@") and then assign the value. Values
# that still contain two dots (e.g. "28.771.43") are never parseable as numbers
# by Excel, so no extra formatting step is required for those.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.771.43"
$ws.Range("E2").Value = "  -2.70%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.886.91"

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.25%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.87"
$ws.Range("E5").Value = "  -1.80%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.23%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4566"
$ws.Range("E7").Value = "  -1.75%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3801"
$ws.Range("E8").Value = "  -3.67%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07715"
$ws.Range("E9").Value = "  -2.71%  "

# Row 10 - Polygon
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9617"
$ws.Range("E10").Value = "  -3.89%  "

# Row 11 - Solana
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.02"
$ws.Range("E11").Value = "  -3.01%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "1.881.19"
$ws.Range("E12").Value = "  -6.89%  "

# Row 13 - Chainlink
$ws.Range("E13").Value = "  -3.68%  "

# Row 14 - Polkadot
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.646"
$ws.Range("E14").Value = "  -3.62%  "

# Row 15 - TRON
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07001"
$ws.Range("E15").Value = "  -1.45%  "

# Row 16 - BinanceUSD
$ws.Range("E16").Value = "  +0.26%  "

# Row 17 - Litecoin
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "83.03"
$ws.Range("E17").Value = "  -6.35%  "

# Row 18 - ShibaInu
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009478"
$ws.Range("E18").Value = "  -5.05%  "

# Row 19 - Avalanche
$ws.Range("E19").Value = "  -3.41%  "

# Row 20 - Dai
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.28%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "28.722.50"
$ws.Range("E21").Value = "  -3.06%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.298"
$ws.Range("E22").Value = "  -4.08%  "

# Row 23 - Cosmos
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.87"
$ws.Range("E23").Value = "  -3.49%  "

# Row 24 - WrappedliquidstakedEther2.0
$ws.Range("D24").Value = "2.128.36"
$ws.Range("E24").Value = "  -5.32%  "

# Row 25 - Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.077"
$ws.Range("E25").Value = "  -2.38%  "

# Row 26 - Monero
$ws.Range("E26").Value = "  -1.27%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.97"
$ws.Range("E27").Value = "  -3.42%  "

# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.586"
$ws.Range("E28").Value = "  -7.04%  "

# Row 29 - BitcoinCash
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.78"
$ws.Range("E29").Value = "  -3.03%  "

# Row 30 - LidoDAOToken
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.809"
$ws.Range("E30").Value = "  -6.63%  "

# Row 31 - Stellar
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09258"
$ws.Range("E31").Value = "  -1.71%  "

# Row 32 - ImmutableX
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8445"
$ws.Range("E32").Value = "  -5.11%  "

# Row 33 - Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.066"
$ws.Range("E33").Value = "  -3.82%  "

# Row 34 - ARBITRUM
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.235"
$ws.Range("E34").Value = "  -8.21%  "

# Row 35 - HuobiToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.008"
$ws.Range("E35").Value = "  -5.13%  "

# Row 36 - Hedera
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05664"
$ws.Range("E36").Value = "  -2.47%  "

# Row 37 - TrustWalletToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.144"
$ws.Range("E37").Value = "  -2.80%  "

# Row 38 - Frax
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.002"
$ws.Range("E38").Value = "  +0.30%  "

# Row 39 - VeChain
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02025"
$ws.Range("E39").Value = "  -4.83%  "

# Row 40 - TheSandbox
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5482"
$ws.Range("E40").Value = "  -4.68%  "

# Row 41 - FraxShare
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.380"
$ws.Range("E41").Value = "  -6.80%  "

# Row 42 - Algorand
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1749"
$ws.Range("E42").Value = "  -4.05%  "

# Row 43 - PEPE
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000002971"
$ws.Range("E43").Value = "  -29.23%  "

# Row 44 - Aptos
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.131"
$ws.Range("E44").Value = "  -7.04%  "

# Row 45 - MXToken
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.697"
$ws.Range("E45").Value = "  +2.08%  "

# Row 46 - EnergySwap/Decentraland swap
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.29"
$ws.Range("E46").Value = "  -6.47%  "

# Row 47 - Decentraland/EnergySwap swap
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5149"
$ws.Range("E47").Value = "  -4.28%  "

# Row 48 - Cronos
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06801"
$ws.Range("E48").Value = "  -2.76%  "

# Row 49 - RenderToken
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.057"
$ws.Range("E49").Value = "  -5.40%  "

# Row 50 - Quant
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.42"
$ws.Range("E50").Value = "  -2.60%  "

# Row 51 - NEARProtocol
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.779"
$ws.Range("E51").Value = "  -5.11%  "
